$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-31: new HKL-reflection data table (2 new entries "Holden" and
# "Rizzie Spiral" inserted after "Spiral5"; "Thomas Hex" renamed to
# "Matthies Hex"; all C:T values recomputed from the rerun simulation).
$rowsData = @(
    @{ row = 4; a = 2; b = 'Holden'; vals = @(0.997192584478829, 1.000701857040393, 1.000701857040393, 0.997192584478829, 1.000701857040393, 1.000701857040393, 0.9984559229869511, 1.001871610715191, 1.000701857040393, 1.000701857040393, 0.9989472207596108, 0.9989472207596108, 0.998783454835391, 0.9995320995198714, 0.9995320995198714, 0.9998245389000018, 0.9998245389000018, 0.9999376148836916) },
    @{ row = 5; a = 3; b = 'Rizzie Spiral'; vals = @(0.994513752543125, 1.001371574925625, 1.001371574925625, 0.994513752543125, 1.001371574925625, 1.001371574925625, 0.996982563572071, 1.003657505162014, 1.001371574925625, 1.001371574925625, 0.9979426637343749, 0.9979426637343749, 0.9976226303469403, 0.9990856341314581, 0.9990856341314581, 0.9996571193299998, 0.9996571193299998, 0.999878091009014) },
    @{ row = 6; a = 4; b = 'RotRing OmegaMax-90'; vals = @(0.9984142526515971, 1.000396441248276, 1.000396441248276, 0.9984142526515971, 1.000396441248276, 1.000396441248276, 0.9991278411719423, 1.001057162099064, 1.000396441248276, 1.000396441248276, 0.9994053469499368, 0.9994053469499368, 0.9993128450239386, 0.9997357117160499, 0.9997357117160499, 0.9999008940991065, 0.9999008940991065, 0.9999647632779055) },
    @{ row = 7; a = 5; b = 'Equal Angle'; vals = @(0.9970521631988478, 1.000736962716139, 1.000736962716139, 0.9970521631988478, 1.000736962716139, 1.000736962716139, 0.9983786912247852, 1.001965226174351, 1.000736962716139, 1.000736962716139, 0.9988945629574935, 0.9988945629574935, 0.9987226057132573, 0.9995086962103753, 0.9995086962103753, 0.9998157628368163, 0.9998157628368163, 0.999934494791067) },
    @{ row = 8; a = 6; b = 'Tilt Rotate'; vals = @(0.9901958024679519, 1.002451052123913, 1.002451052123913, 0.9901958024679519, 1.002451052123913, 1.002451052123913, 0.9946076941840079, 1.006536151903368, 1.002451052123913, 1.002451052123913, 0.9963234272959325, 0.9963234272959325, 0.9957515162586242, 0.9983659689052594, 0.9983659689052594, 0.9993872397099228, 0.9993872397099228, 0.9997821341545112) },
    @{ row = 9; a = 7; b = 'CLR'; vals = @(0.9996751277528185, 1.000081222221127, 1.000081222221127, 0.9996751277528185, 1.000081222221127, 1.000081222221127, 0.9998213225223261, 1.000216580719323, 1.000081222221127, 1.000081222221127, 0.9998781749869725, 0.9998781749869725, 0.9998592241654237, 0.9999458573983572, 0.9999458573983572, 0.9999796986040496, 0.9999796986040496, 0.9999927829429746) },
    @{ row = 10; a = 8; b = 'Rizzie Hex'; vals = @(0.9999755728388378, 1.000006110863086, 1.000006110863086, 0.9999755728388378, 1.000006110863086, 1.000006110863086, 0.9999865667540628, 1.000016284515201, 1.000006110863086, 1.000006110863086, 0.9999908418509618, 0.9999908418509618, 0.9999894168186622, 0.9999959315216698, 0.9999959315216698, 0.9999984763570238, 0.9999984763570238, 0.9999994594495599) },
    @{ row = 11; a = 9; b = 'Matthies Hex'; vals = @(0.9994607561731895, 1.000134817007111, 1.000134817007111, 0.9994607561731895, 1.000134817007111, 1.000134817007111, 0.999703416915124, 1.000359495831011, 1.000134817007111, 1.000134817007111, 0.9997977865901504, 0.9997977865901504, 0.9997663300318083, 0.9999101300624705, 0.9999101300624705, 0.9999663017986307, 0.9999663017986307, 0.9999880199901097) },
    @{ row = 12; a = 10; b = 'Tilt Rotate_Partial'; vals = @(0.9900493723235453, 1.00248766002167, 1.00248766002167, 0.9900493723235453, 1.00248766002167, 1.00248766002167, 0.994527157302363, 1.006633772398039, 1.00248766002167, 1.00248766002167, 0.9962685161726076, 0.9962685161726076, 0.9956880632158595, 0.998341564122295, 0.998341564122295, 0.9993780880971388, 0.9993780880971388, 0.9997788803481593) },
    @{ row = 13; a = 11; b = 'RotRing OmegaMax-60'; vals = @(0.9893081457438621, 1.002672968615168, 1.002672968615168, 0.9893081457438621, 1.002672968615168, 1.002672968615168, 0.994119481923199, 1.007127898833531, 1.002672968615168, 1.002672968615168, 0.9959905571795149, 0.9959905571795149, 0.9953668654274096, 0.998218027658066, 0.998218027658066, 0.9993317628973415, 0.9993317628973415, 0.999762405391016) },
    @{ row = 14; a = 12; b = 'Equal Angle_Partial'; vals = @(0.9966563488210515, 1.000835916273683, 1.000835916273683, 0.9966563488210515, 1.000835916273683, 1.000835916273683, 0.9981609925473705, 1.002229103852631, 1.000835916273683, 1.000835916273683, 0.9987461325473672, 0.9987461325473672, 0.9985510858807016, 0.9994427271228058, 0.9994427271228058, 0.9997910244105251, 0.9997910244105251, 0.999925699007017) },
    @{ row = 15; a = 13; b = 'Rizzie Hex_Partial'; vals = @(1.001957007022278, 0.9995107532491248, 0.9995107532491248, 1.001957007022278, 0.9995107532491248, 0.9995107532491248, 1.001076354526652, 0.9986953291366857, 0.9995107532491248, 0.9995107532491248, 1.000733880135701, 1.000733880135701, 1.000848038266018, 1.000326171173509, 1.000326171173509, 1.000122316692413, 1.000122316692413, 1.000043491738832) },
    @{ row = 16; a = 14; b = 'ND Single'; vals = @(0.9828084099999992, 1.004297899999999, 1.004297899999999, 0.9828084099999992, 1.004297899999999, 1.004297899999999, 0.9905446300000009, 1.011461099999999, 1.004297899999999, 1.004297899999999, 0.9935531549999991, 0.9935531549999991, 0.9925503133333331, 0.9971347366666657, 0.9971347366666657, 0.9989255274999991, 0.9989255274999991, 0.9996179733333328) },
    @{ row = 17; a = 15; b = 'RD Single'; vals = @(1.0687664, 0.98280841, 0.98280841, 1.0687664, 0.98280841, 0.98280841, 1.0378215, 0.95415576, 0.98280841, 0.98280841, 1.025787405, 1.025787405, 1.02979877, 1.011461073333333, 1.011461073333333, 1.0042979075, 1.0042979075, 1.001528148333333) },
    @{ row = 18; a = 16; b = 'TD Single'; vals = @(0.98280841, 1.0042979, 1.0042979, 0.98280841, 1.0042979, 1.0042979, 0.99054463, 1.0114611, 1.0042979, 1.0042979, 0.9935531550000001, 0.9935531550000001, 0.9925503133333334, 0.9971347366666667, 0.9971347366666667, 0.9989255275000001, 0.9989255275000001, 0.9996179733333334) },
    @{ row = 19; a = 17; b = 'Morris Single'; vals = @(1.0472769, 0.98818078, 0.98818078, 1.0472769, 0.98818078, 0.98818078, 1.0260023, 0.96848208, 0.98818078, 0.98818078, 1.01772884, 1.01772884, 1.02048666, 1.007879486666667, 1.007879486666667, 1.00295481, 1.00295481, 1.001050603333334) },
    @{ row = 20; a = 18; b = 'Ring Perpendicular to ND'; vals = @(0.9945834764383562, 1.001354144520548, 1.001354144520548, 0.9945834764383562, 1.001354144520548, 1.001354144520548, 0.9970209117808218, 1.003611022465753, 1.001354144520548, 1.001354144520548, 0.997968810479452, 0.997968810479452, 0.9976528442465753, 0.9990972551598173, 0.9990972551598173, 0.9996614774999999, 0.9996614774999999, 0.9998796407077624) },
    @{ row = 21; a = 19; b = 'Ring Perpendicular to RD'; vals = @(1.023525342105263, 0.9941186657894736, 0.9941186657894736, 1.023525342105263, 0.9941186657894736, 0.9941186657894736, 1.012938931578947, 0.984316445263158, 0.9941186657894736, 0.9941186657894736, 1.008822003947368, 1.008822003947368, 1.010194313157895, 1.00392089122807, 1.00392089122807, 1.001470334868421, 1.001470334868421, 1.000522786052631) },
    @{ row = 22; a = 20; b = 'Ring Perpendicular to TD'; vals = @(0.9952496973684208, 1.001187588947368, 1.001187588947368, 0.9952496973684208, 1.001187588947368, 1.001187588947368, 0.9973873331578945, 1.003166876315789, 1.001187588947368, 1.001187588947368, 0.9982186431578944, 0.9982186431578944, 0.9979415398245611, 0.9992082917543857, 0.9992082917543857, 0.9997031160526313, 0.9997031160526313, 0.9998944456140348) },
    @{ row = 23; a = 21; b = 'OffsetFTD'; vals = @(1.045701945667218, 0.988574515860675, 0.988574515860675, 1.045701945667218, 0.988574515860675, 0.988574515860675, 1.025136069581616, 0.969532040033978, 0.988574515860675, 0.988574515860675, 1.017138230763946, 1.017138230763946, 1.019804177036503, 1.007616992462856, 1.007616992462856, 1.002856373312311, 1.002856373312311, 1.001015600477473) },
    @{ row = 24; a = 22; b = 'OffsetATD'; vals = @(1.012181208237609, 0.9969546944137678, 0.9969546944137678, 1.012181208237609, 0.9969546944137678, 0.9969546944137678, 1.006699670278214, 0.9918791945001935, 0.9969546944137678, 0.9969546944137678, 1.004567951325689, 1.004567951325689, 1.005278524309864, 1.002030199021715, 1.002030199021715, 1.000761322869728, 1.000761322869728, 1.000270692709553) },
    @{ row = 25; a = 23; b = 'OffsetF45'; vals = @(0.9995196357321918, 1.000120094216192, 1.000120094216192, 0.9995196357321918, 1.000120094216192, 1.000120094216192, 0.9997357961980131, 1.000320241021014, 1.000120094216192, 1.000120094216192, 0.9998198649741917, 0.9998198649741917, 0.9997918420487988, 0.9999199413881916, 0.9999199413881916, 0.9999699795951916, 0.9999699795951916, 0.9999893259332989) },
    @{ row = 26; a = 24; b = 'OffsetA45'; vals = @(0.9998696851522243, 1.000032580109521, 1.000032580109521, 0.9998696851522243, 1.000032580109521, 1.000032580109521, 0.9999283311509329, 1.000086877963037, 1.000032580109521, 1.000032580109521, 0.9999511326308728, 0.9999511326308728, 0.9999435321375595, 0.9999782817904223, 0.9999782817904223, 0.9999918563701971, 0.9999918563701971, 0.9999971057657929) },
    @{ row = 27; a = 25; b = 'OffsetFRD'; vals = @(0.9642704642877111, 1.008932394824936, 1.008932394824936, 0.9642704642877111, 1.008932394824936, 1.008932394824936, 0.9803487537850353, 1.023819690784401, 1.008932394824936, 1.008932394824936, 0.9866014295563237, 0.9866014295563237, 0.9845172042992276, 0.9940450846458612, 0.9940450846458612, 0.9977669121906299, 0.9977669121906299, 0.9992060155553261) },
    @{ row = 28; a = 26; b = 'OffsetARD'; vals = @(0.9904061168911402, 1.002398478429737, 1.002398478429737, 0.9904061168911402, 1.002398478429737, 1.002398478429737, 0.9947233633050211, 1.006395918014905, 1.002398478429737, 1.002398478429737, 0.9964022976604388, 0.9964022976604388, 0.9958426528752996, 0.9984010245835382, 0.9984010245835382, 0.9994003880450879, 0.9994003880450879, 0.9997868055833797) },
    @{ row = 29; a = 27; b = 'Gaussian Quadrature'; vals = @(0.9989753889113621, 1.000256170142117, 1.000256170142117, 0.9989753889113621, 1.000256170142117, 1.000256170142117, 0.9994364608848769, 1.000683081276603, 1.000256170142117, 1.000256170142117, 0.9996157795267397, 0.9996157795267397, 0.9995560066461188, 0.9998292430651988, 0.9998292430651988, 0.9999359748344285, 0.9999359748344285, 0.9999772402498657) },
    @{ row = 30; a = 28; b = 'Michael-CCHex'; vals = @(1.001033708190077, 0.9997415824313458, 0.9997415824313458, 1.001033708190077, 0.9997415824313458, 0.9997415824313458, 1.000568549101501, 0.9993108581496504, 0.9997415824313458, 0.9997415824313458, 1.000387645310711, 1.000387645310711, 1.000447946574308, 1.00017229101759, 1.00017229101759, 1.000064613871029, 1.000064613871029, 1.000022977122544) },
    @{ row = 31; a = 29; b = 'Michael-SNHex'; vals = @(1.004495382678261, 0.998876172065593, 0.998876172065593, 1.004495382678261, 0.998876172065593, 0.998876172065593, 1.002472472718944, 0.9970030735112908, 0.998876172065593, 0.998876172065593, 1.001685777371927, 1.001685777371927, 1.001948009154266, 1.000749242269816, 1.000749242269816, 1.00028097471876, 1.00028097471876, 1.000099907517546) }
)

foreach ($r in $rowsData) {
    $rowNum = $r.row
    $ws.Cells.Item($rowNum, 1).Value2 = $r.a
    $ws.Cells.Item($rowNum, 2).Value2 = $r.b
    $vals = $r.vals
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($rowNum, $i + 3).Value2 = $vals[$i]
    }
}

# Column A cells use the bold/bordered header style (style index 1 in the
# original workbook) for every data row, same as rows 2-29 already had.
$ws.Range("A4:A5").Style = $ws.Range("A3").Style
$ws.Range("A30:A31").Style = $ws.Range("A29").Style

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()